$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: shift the header number sequence ---
$ws.Range("E1").ClearContents()
$ws.Range("G1").Value = 18
$ws.Range("H1").Value = 19
$ws.Range("I1").Value = 20
$ws.Range("J1").Value = 21
$ws.Range("K1").Value = 22
$ws.Range("L1").Value = 23
$ws.Range("M1").Value = 24
$ws.Range("N1").Value = 25
$ws.Range("O1").Value = 26
$ws.Range("P1").Value = 27
$ws.Range("Q1").Value = 28

# --- Row 5 / Row 6: new "OK" markers ---
$ws.Range("F5").Value = "OK"
$ws.Range("G6").Value = "OK"

# --- Row 7: the "x" marker moves from G7 to H7 ---
$ws.Range("G7").ClearContents()
$ws.Range("H7").Value = "x"

# --- Row 9 (new row) ---
$ws.Range("C9").Value = 9
$ws.Range("I9").Value = "x"

# --- Row 10: drop J10, add N10:P10 ---
$ws.Range("J10").ClearContents()
$ws.Range("N10").Value = "x"
$ws.Range("O10").Value = "x"
$ws.Range("P10").Value = "x"

# --- Row 11 (new row) ---
$ws.Range("C11").Value = "Rette"

# --- Shade the L:M helper column block (rows 4-10) ---
$ws.Range("L4:M10").Interior.Color = 15921906

# --- Selection moves to I9 ---
$ws.Range("I9").Select() | Out-Null
